# Update "Lũy kế tháng HỆ THỐNG" report data (Notion export) for rows 7-12:
# refresh last_edited_time and the recalculated formula/number properties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newEditedTime = "2024-07-24T16:01:00.000Z"

# Row 7 (Thang 7)
$ws.Range("D7").Value = $newEditedTime
$ws.Range("T7").Value = 56300000
$ws.Range("W7").Value = 267881000
$ws.Range("AA7").Value = 256587000
$ws.Range("AE7").Value = 524468000
$ws.Range("AH7").Value = 449968000
$ws.Range("AK7").Value = 70
$ws.Range("AQ7").Value = 506268000

# Row 8 (Thang 6)
$ws.Range("D8").Value = $newEditedTime
$ws.Range("T8").Value = 118600000
$ws.Range("AA8").Value = 748973000
$ws.Range("AE8").Value = 1080200000
$ws.Range("AH8").Value = 931300000

# Row 9 (Thang 5)
$ws.Range("D9").Value = $newEditedTime
$ws.Range("T9").Value = 133786000
$ws.Range("AA9").Value = 407109000
$ws.Range("AE9").Value = 439686000
$ws.Range("AH9").Value = 368300000

# Row 10 (Thang 4)
$ws.Range("D10").Value = $newEditedTime
$ws.Range("T10").Value = 232400000
$ws.Range("AA10").Value = 342764000
$ws.Range("AE10").Value = 842750000
$ws.Range("AH10").Value = 783350000

# Row 11 (Thang 3)
$ws.Range("D11").Value = $newEditedTime
$ws.Range("T11").Value = 166100000
$ws.Range("AA11").Value = 455643000
$ws.Range("AE11").Value = 779450000
$ws.Range("AH11").Value = 723950000

# Row 12 (Thang 2)
$ws.Range("D12").Value = $newEditedTime
$ws.Range("T12").Value = 49498000
$ws.Range("AA12").Value = 261993000
$ws.Range("AE12").Value = 377795000
$ws.Range("AH12").Value = 349795000
